$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6712798476219177
$ws.Range("B1").Value = 1.462276935577393
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.331178665161133
$ws.Range("E1").Value = 1.386141538619995
